$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Natmi re-run following Dr Hou advice: for every Gas6-Mertk sending/target-cluster
# pair (rows 2-17), the ligand- and receptor-expressing cell counts (columns E, K)
# increase from 1 to 3, and every statistic derived from them (average/total
# expression values, derived-specificity scores, and edge weights) is refreshed
# to the recomputed NATMI output.

$columns = @("E", "G", "H", "I", "J", "K", "M", "N", "O", "P", "Q", "R", "S", "T")

$newValues = @{
    2 = @(3, 22.06580366666666, 66.19741099999999, 0.1403713626377477, 0.1403713626377477, 3, 11.76011933333333, 35.280358, 0.1622550713151862, 0.1622550713151862, 259.4964843059042, 2335.468358753138, 0.02277596545539762, 0.02277596545539762)
    3 = @(3, 22.06580366666666, 66.19741099999999, 0.1403713626377477, 0.1403713626377477, 3, 3.903127666666666, 11.709383, 0.05385168636105758, 0.05385168636105758, 86.12564877860143, 775.1308390074128, 0.007559234594842267, 0.007559234594842266)
    4 = @(3, 22.06580366666666, 66.19741099999999, 0.1403713626377477, 0.1403713626377477, 3, 50.38343933333334, 151.150318, 0.6951433323438234, 0.6951433323438234, 1111.751080380744, 10005.7597234267, 0.09757821678964723, 0.09757821678964722)
    5 = @(3, 22.06580366666666, 66.19741099999999, 0.1403713626377477, 0.1403713626377477, 3, 6.432523333333333, 19.29757, 0.08874990997993266, 0.08874990997993268, 141.9387969545855, 1277.44917259127, 0.0124579457978606, 0.0124579457978606)
    6 = @(3, 98.09611, 294.28833, 0.6240373039738243, 0.6240373039738243, 3, 11.76011933333333, 35.280358, 0.1622550713151862, 0.1622550713151862, 1153.621959735793, 10382.59763762214, 0.1012532172596094, 0.1012532172596094)
    7 = @(3, 98.09611, 294.28833, 0.6240373039738243, 0.6240373039738243, 3, 3.903127666666666, 11.709383, 0.05385168636105758, 0.05385168636105758, 382.8816409333766, 3445.934768400389, 0.03360546117119834, 0.03360546117119834)
    8 = @(3, 98.09611, 294.28833, 0.6240373039738243, 0.6240373039738243, 3, 50.38343933333334, 151.150318, 0.6951433323438234, 0.6951433323438234, 4942.419407020993, 44481.77466318894, 0.4337953709912197, 0.4337953709912197)
    9 = @(3, 98.09611, 294.28833, 0.6240373039738243, 0.6240373039738243, 3, 6.432523333333333, 19.29757, 0.08874990997993266, 0.08874990997993268, 631.0055164842333, 5679.049648358099, 0.05538325455179678, 0.05538325455179679)
    10 = @(3, 29.56610533333334, 88.69831600000001, 0.1880844476016372, 0.1880844476016372, 3, 11.76011933333333, 35.280358, 0.1622550713151862, 0.1622550713151862, 347.7009269419032, 3129.308342477128, 0.03051765545888105, 0.03051765545888105)
    11 = @(3, 29.56610533333334, 88.69831600000001, 0.1880844476016372, 0.1880844476016372, 3, 3.903127666666666, 11.709383, 0.05385168636105758, 0.05385168636105758, 115.4002837221142, 1038.602553499028, 0.01012866468163614, 0.01012866468163614)
    12 = @(3, 29.56610533333334, 88.69831600000001, 0.1880844476016372, 0.1880844476016372, 3, 50.38343933333334, 151.150318, 0.6951433323438234, 0.6951433323438234, 1489.642074384943, 13406.77866946449, 0.1307456496678494, 0.1307456496678494)
    13 = @(3, 29.56610533333334, 88.69831600000001, 0.1880844476016372, 0.1880844476016372, 3, 6.432523333333333, 19.29757, 0.08874990997993266, 0.08874990997993268, 190.1846624324578, 1711.66196189212, 0.01669247779327067, 0.01669247779327067)
    14 = @(3, 7.467887999999999, 22.403664, 0.04750688578679087, 0.04750688578679088, 3, 11.76011933333333, 35.280358, 0.1622550713151862, 0.1622550713151862, 87.82325404796799, 790.409286431712, 0.007708233141298158, 0.007708233141298159)
    15 = @(3, 7.467887999999999, 22.403664, 0.04750688578679087, 0.04750688578679088, 3, 3.903127666666666, 11.709383, 0.05385168636105758, 0.05385168636105758, 29.148120264368, 262.333082379312, 0.002558325913380846, 0.002558325913380847)
    16 = @(3, 7.467887999999999, 22.403664, 0.04750688578679087, 0.04750688578679088, 3, 50.38343933333334, 151.150318, 0.6951433323438234, 0.6951433323438234, 376.257881996128, 3386.320937965152, 0.03302409489510723, 0.03302409489510724)
    17 = @(3, 7.467887999999999, 22.403664, 0.04750688578679087, 0.04750688578679088, 3, 6.432523333333333, 19.29757, 0.08874990997993266, 0.08874990997993268, 48.03736381072, 432.33627429648, 0.004216231837004632, 0.004216231837004634)
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
